$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2 through 39 from serial date 45182 to 45184
$ws.Range("C2:C39").Value = 45184
